$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 (Week1 "Course" row) -------------------------------------
# Add the course times for Wed/Thu/Fri that were previously blank.
$ws.Range("C13").Value = "4->5"
$ws.Range("D13").Value = "6->7"
$ws.Range("F13").Value = "3->5"

# E13 stays empty but gains a bottom border (new style): reset it to the
# plain "Normal" style first so no stray number-format/font carries over,
# then draw just the bottom border.
$ws.Range("E13").Style = "Normal"
$ws.Range("E13").Borders.Item(9).LineStyle = 1

# --- Row 16: a third meeting label, bold like the others -------------
$ws.Range("F16").Value = "Meeting 3"
$ws.Range("F16").Font.Bold = $true

# --- Week2 meetings section ------------------------------------------
$ws.Range("D33").Value = "Meeting 4"
$ws.Range("D33").Font.Bold = $true

$ws.Range("D31").Value = "10.30 -> Craig, Lorna"

# --- Selection / view state -------------------------------------------
# Scroll the window so row 18 is the top visible row, then select E34
# (matches the saved sheetView's topLeftCell/selection in the workbook).
try {
    $excel.ActiveWindow.ScrollRow = 18
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # view-scroll state isn't critical to the data edit; ignore if unsupported
}
$ws.Range("E34").Select()
